$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A64").Value = "2025/12/05 09:00"
$ws.Range("B64").Value = "-"
$ws.Range("C64").Value = "-"
$ws.Range("D64").Value = "-"
$ws.Range("E64").Value = "-"
$ws.Range("F64").Value = "-"
$ws.Range("G64").Value = "-"
